$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new localization key ("strWindowPos") was added to the source project, for
# the new "Remember window position and size on startup" setting, shown on
# the "User interface" tab of the "settings" form. This inserts a new table
# row for it, right before the current row 32 ("strDlgReset"), pushing the
# rest of the translation table down by one row.
$ws.Range("B32:F32").Insert()

# Keep the "Tabla13" table definition (and its filter) in sync with the newly
# inserted row, growing it by exactly one row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))

# Fill in the new row: File / Key / Comment / English columns (the
# "Hungarian (hu-HU)" translation column F is intentionally left blank, same
# as other not-yet-translated rows).
$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = "In ""settings"" form, tab ""User interface"""
$ws.Range("E32").Value = "Remember window position and size on startup"

# The pre-existing "strChkDlgPath" row (row 25) belongs to the same
# "settings" form / "User interface" tab, so it now gets the same
# clarifying comment.
$ws.Range("D25").Value = "In ""settings"" form, tab ""User interface"""

# The new, longer comment text no longer fits the old "Comment" column
# width, so widen it to fit (mirrors the author widening column D by hand).
$ws.Columns.Item(4).ColumnWidth = 34.91
